# LoginData.xlsx update: add "register" and "tryEditor" worksheets with
# registration / python-editor test data, and touch up styling / selection
# on the existing "login" sheet.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Existing "login" sheet: move the active selection, drop tab-select.
# ---------------------------------------------------------------------
$login = $wb.Worksheets.Item("login")
$login.Range("B5").Select()

# ---------------------------------------------------------------------
# 2. New "register" worksheet (placed right after "login").
# ---------------------------------------------------------------------
$register = $wb.Worksheets.Add($null, $login)
$register.Name = "register"

# Seed every used cell with the existing thin-all-round border style
# (copied from the already-styled "login" sheet so the stylesheet entry
# is reused rather than rebuilt from scratch).
$login.Range("A1").Copy()
$register.Range("A1:D10").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$register.Range("A1").Value = "username"
$register.Range("B1").Value = "password"
$register.Range("C1").Value = "password confirmation"
$register.Range("D1").Value = "message"

$register.Range("A2").Value = "techtitanstest4"
$register.Range("B2").Value = "Time2fly$"
$register.Range("C2").Value = "Time2fly$"
$register.Range("D2").Value = "New Account Created. You are logged in as"

$register.Range("D3").Value = "Please fill out this field."

$register.Range("A4").Value = "techtitanstest"
$register.Range("D4").Value = "Please fill out this field."

$register.Range("A5").Value = "techtitanstest"
$register.Range("B5").Value = "Time2fly$"
$register.Range("D5").Value = "Please fill out this field."

$register.Range("A6").Value = "#"
$register.Range("B6").Value = "Time2fly$"
$register.Range("C6").Value = "Time2fly$"
$register.Range("D6").Value = "password_mismatch:The two password fields didn’t match."

$register.Range("A7").Value = "techtitanstest"
$register.Range("B7").Value = "Time2fly$"
$register.Range("C7").Value = "Time2fly"
$register.Range("D7").Value = "password_mismatch:The two password fields didn’t match."

$register.Range("A8").Value = "techtitanstest"
$register.Range("B8").Value = "Time"
$register.Range("C8").Value = "Time"
$register.Range("D8").Value = "password_mismatch:The two password fields didn’t match."

$register.Range("A9").Value = "techtitanstest"
$register.Range("B9").Value = 12345678901
$register.Range("C9").Value = 12345678901
$register.Range("D9").Value = "password_mismatch:The two password fields didn’t match."

# D1 only keeps its left/right edges (no top/bottom rule).
$register.Range("D1").Borders(8).LineStyle = 0
$register.Range("D1").Borders(9).LineStyle = 0

# Rows 2-9 additionally get a Text number format, applied after the
# values so the two numeric-looking password cells (B9/C9) keep their
# numeric storage instead of turning into text.
$register.Range("A2:D9").NumberFormat = "@"

$register.Range("A2").Select()

# ---------------------------------------------------------------------
# 3. New "tryEditor" worksheet (placed right after "register").
# ---------------------------------------------------------------------
$tryEditor = $wb.Worksheets.Add($null, $register)
$tryEditor.Name = "tryEditor"

$login.Range("A1").Copy()
$tryEditor.Range("A1:D7").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$tryEditor.Range("A1").Value = "username"
$tryEditor.Range("B1").Value = "password"
$tryEditor.Range("C1").Value = "python code"
$tryEditor.Range("D1").Value = "alert message"

$tryEditor.Range("A2").Value = "techtitans"
$tryEditor.Range("B2").Value = "Time2fly$"

$tryEditor.Range("A3").Value = "techtitans"
$tryEditor.Range("B3").Value = "Time2fly$"

$tryEditor.Range("A4").Value = "techtitans"
$tryEditor.Range("B4").Value = "Time2fly$"

$tryEditor.Range("A5").Value = "techtitans"
$tryEditor.Range("B5").Value = "Time2fly$"

$tryEditor.Range("A6").Value = "techtitans"
$tryEditor.Range("B6").Value = "Time2fly$"
$tryEditor.Range("C6").Value = 'print"hello"'

$tryEditor.Range("A7").Value = "techtitans"
$tryEditor.Range("B7").Value = "Time2fly$"
$tryEditor.Range("C7").Value = 'print"hello";;'
$tryEditor.Range("D7").Value = "SyntaxError: bad input on line 1"

$tryEditor.Range("D7").Select()

# ---------------------------------------------------------------------
# 4. Make "register" the active sheet/tab, matching the saved view state.
# ---------------------------------------------------------------------
$register.Activate()
